$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Halwan, Sharjah'
$ws.Range("B2").Value = '100,000 AED/year'
$ws.Range("C2").Value = '12,000 sqft'
$ws.Range("D2").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-halwan-13132606.html'

$ws.Range("A3").Value = 'Al Fisht, Al Heerah, Sharjah'
$ws.Range("B3").Value = '110,000 AED/year'
$ws.Range("C3").Value = '10,115 sqft'
$ws.Range("D3").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-al-heerah-al-fisht-12934465.html'

$ws.Range("A4").Value = 'Amargo, Damac Hills 2, Dubai'
$ws.Range("B4").Value = '120,000 AED/year'
$ws.Range("C4").Value = '1,208 sqft'
$ws.Range("D4").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-damac-hills-2-amargo-13137496.html'

$ws.Range("A5").Value = 'Mirdif Villas, Mirdif, Dubai'
$ws.Range("B5").Value = '145,000 AED/year'
$ws.Range("C5").Value = '10,125 sqft'
$ws.Range("D5").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-mirdif-mirdif-villas-12805233.html'

$ws.Range("A6").Value = 'Hoshi, Al Badie, Sharjah'
$ws.Range("B6").Value = '160,000 AED/year'
$ws.Range("C6").Value = '8,500 sqft'
$ws.Range("D6").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-al-badie-hoshi-13144570.html'

$ws.Range("A7").Value = 'Al Rahmaniya, Sharjah'
$ws.Range("B7").Value = '160,000 AED/year'
$ws.Range("C7").Value = '20,000 sqft'
$ws.Range("D7").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-al-rahmaniya-13132605.html'

$ws.Range("A8").Value = 'Sharqan, Al Heerah, Sharjah'
$ws.Range("B8").Value = '160,000 AED/year'
$ws.Range("C8").Value = '8,000 sqft'
$ws.Range("D8").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-al-heerah-sharqan-13082488.html'

$ws.Range("A9").Value = 'Al Fisht, Al Heerah, Sharjah'
$ws.Range("B9").Value = '165,000 AED/year'
$ws.Range("C9").Value = '2,750 sqft'
$ws.Range("D9").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-al-heerah-al-fisht-13127318.html'

$ws.Range("A10").Value = 'Abu Hail, Deira, Dubai'
$ws.Range("B10").Value = '180,000 AED/year'
$ws.Range("C10").Value = '3,916 sqft'
$ws.Range("D10").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-deira-abu-hail-12923594.html'

$ws.Range("A11").Value = 'Hoshi, Al Badie, Sharjah'
$ws.Range("B11").Value = '180,000 AED/year'
$ws.Range("C11").Value = '6,000 sqft'
$ws.Range("D11").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-al-badie-hoshi-13131146.html'

$ws.Range("A12").Value = 'Hoshi, Al Badie, Sharjah'
$ws.Range("B12").Value = '185,000 AED/year'
$ws.Range("C12").Value = '9,000 sqft'
$ws.Range("D12").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-al-badie-hoshi-13135794.html'

$ws.Range("A13").Value = 'Hoshi 1, Hoshi, Al Badie, Sharjah'
$ws.Range("B13").Value = '185,000 AED/year'
$ws.Range("C13").Value = '9,200 sqft'
$ws.Range("D13").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-sharjah-al-badie-hoshi-hoshi-1-13136019.html'

$ws.Range("A14").Value = 'Al Aweer 1, Al Aweer, Dubai'
$ws.Range("B14").Value = '189,999 AED/year'
$ws.Range("C14").Value = '12,056 sqft'
$ws.Range("D14").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-al-aweer-al-aweer-1-13165462.html'

$ws.Range("A15").Value = 'Al Rashidiya, Dubai'
$ws.Range("B15").Value = '190,000 AED/year'
$ws.Range("C15").Value = '4,000 sqft'
$ws.Range("D15").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-al-rashidiya-13033334.html'

$ws.Range("A16").Value = 'Amaranta B, Villanova, Dubai Land, Dubai'
$ws.Range("B16").Value = '199,999 AED/year'
$ws.Range("C16").Value = '2,475 sqft'
$ws.Range("D16").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-dubai-land-villanova-amaranta-b-13164457.html'

$ws.Range("A17").Value = 'Wadi Al Safa 5, Dubai'
$ws.Range("B17").Value = '200,000 AED/year'
$ws.Range("C17").Value = '2,475 sqft'
$ws.Range("D17").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-wadi-al-safa-5-13138917.html'

$ws.Range("A18").Value = '180 Degrees Villas, Liwan, Dubai Land, Dubai'
$ws.Range("B18").Value = '200,000 AED/year'
$ws.Range("C18").Value = '2,684 sqft'
$ws.Range("D18").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-dubai-land-liwan-180-degrees-villas-13074818.html'

$ws.Range("A19").Value = '180 Degrees Villas, Liwan, Dubai Land, Dubai'
$ws.Range("B19").Value = '200,000 AED/year'
$ws.Range("C19").Value = '2,684 sqft'
$ws.Range("D19").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-dubai-land-liwan-180-degrees-villas-13095855.html'

$ws.Range("A20").Value = '180 Degrees Villas, Liwan, Dubai Land, Dubai'
$ws.Range("B20").Value = '200,000 AED/year'
$ws.Range("C20").Value = '2,675 sqft'
$ws.Range("D20").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-dubai-land-liwan-180-degrees-villas-13052035.html'

$ws.Range("A21").Value = 'Al Aweer 1, Al Aweer, Dubai'
$ws.Range("B21").Value = '200,000 AED/year'
$ws.Range("C21").Value = '12,700 sqft'
$ws.Range("D21").Value = 'https://www.propertyfinder.ae/en/plp/rent/villa-for-rent-dubai-al-aweer-al-aweer-1-13145382.html'

$ws.Range("A22:D25").Delete()
